$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.563671284420025
$ws.Cells.Item(2, 3).Value = 0.1639436282541027
$ws.Cells.Item(2, 4).Value = 0.1567694438889262
$ws.Cells.Item(2, 6).Value = 1.627944705169611
$ws.Cells.Item(2, 7).Value = 0.002473875492890764
$ws.Cells.Item(2, 10).Value = 0.1893502527895983
$ws.Cells.Item(2, 12).Value = 0.3784317300754338
$ws.Cells.Item(2, 14).Value = 1.425601805552539
$ws.Cells.Item(2, 15).Value = 4.074241858967383
$ws.Cells.Item(3, 2).Value = 1.466737609897677
$ws.Cells.Item(3, 3).Value = 0.1537511046946065
$ws.Cells.Item(3, 4).Value = 0.1561867015483358
$ws.Cells.Item(3, 6).Value = 1.629938341303188
$ws.Cells.Item(3, 7).Value = 0.002476927565470671
$ws.Cells.Item(3, 10).Value = 0.1909406558013327
$ws.Cells.Item(3, 12).Value = 0.3710209256508676
$ws.Cells.Item(3, 14).Value = 1.438651867270686
$ws.Cells.Item(3, 15).Value = 4.081803531834879
$ws.Cells.Item(4, 2).Value = 1.407639608704415
$ws.Cells.Item(4, 3).Value = 0.1474472678269336
$ws.Cells.Item(4, 4).Value = 0.1558741488387909
$ws.Cells.Item(4, 6).Value = 1.632048475836669
$ws.Cells.Item(4, 7).Value = 0.002478903121655551
$ws.Cells.Item(4, 10).Value = 0.1919847313962535
$ws.Cells.Item(4, 12).Value = 0.3666133597237859
$ws.Cells.Item(4, 14).Value = 1.447208065438218
$ws.Cells.Item(4, 15).Value = 4.088813334630061
$ws.Cells.Item(5, 2).Value = 1.383663772141688
$ws.Cells.Item(5, 3).Value = 0.1448670838160382
$ws.Cells.Item(5, 4).Value = 0.1557582060865599
$ws.Cells.Item(5, 6).Value = 1.633131185704322
$ws.Cells.Item(5, 7).Value = 0.00247973379460753
$ws.Cells.Item(5, 10).Value = 0.192427202990423
$ws.Cells.Item(5, 12).Value = 0.3648532894745813
$ws.Cells.Item(5, 14).Value = 1.450831516309322
$ws.Cells.Item(5, 15).Value = 4.092264828195738
$ws.Cells.Item(6, 2).Value = 1.37968911576462
$ws.Cells.Item(6, 3).Value = 0.1444379668092068
$ws.Cells.Item(6, 4).Value = 0.1557396452686817
$ws.Cells.Item(6, 6).Value = 1.633324426610351
$ws.Cells.Item(6, 7).Value = 0.002479873276836672
$ws.Cells.Item(6, 10).Value = 0.1925017023541695
$ws.Cells.Item(6, 12).Value = 0.3645632132982115
$ws.Cells.Item(6, 14).Value = 1.451441448387726
$ws.Cells.Item(6, 15).Value = 4.092873872301993
$ws.Cells.Item(7, 2).Value = 1.407315826574575
$ws.Cells.Item(7, 3).Value = 0.1474125162020528
$ws.Cells.Item(7, 4).Value = 0.1558725388735027
$ws.Cells.Item(7, 6).Value = 1.632062175464995
$ws.Cells.Item(7, 7).Value = 0.002478914220614771
$ws.Cells.Item(7, 10).Value = 0.1919906298633602
$ws.Cells.Item(7, 12).Value = 0.3665894766201347
$ws.Cells.Item(7, 14).Value = 1.447256378831113
$ws.Cells.Item(7, 15).Value = 4.088857474101189
$ws.Cells.Item(8, 2).Value = 1.53016234831955
$ws.Cells.Item(8, 3).Value = 0.1604387907389366
$ws.Cells.Item(8, 4).Value = 0.1565591489092313
$ws.Cells.Item(8, 6).Value = 1.628448217893585
$ws.Cells.Item(8, 7).Value = 0.002474906812594662
$ws.Cells.Item(8, 10).Value = 0.1898846094210569
$ws.Cells.Item(8, 12).Value = 0.3758469642448716
$ws.Cells.Item(8, 14).Value = 1.429988733989084
$ws.Cells.Item(8, 15).Value = 4.076357720884346
$ws.Cells.Item(9, 2).Value = 1.774338678260847
$ws.Cells.Item(9, 3).Value = 0.1856163979789471
$ws.Cells.Item(9, 4).Value = 0.1582628260721108
$ws.Cells.Item(9, 6).Value = 1.628392956809719
$ws.Cells.Item(9, 7).Value = 0.002467850711595285
$ws.Cells.Item(9, 10).Value = 0.1862901649098454
$ws.Cells.Item(9, 12).Value = 0.3951270654434609
$ws.Cells.Item(9, 14).Value = 1.400434784739446
$ws.Cells.Item(9, 15).Value = 4.070641485079278
$ws.Cells.Item(10, 2).Value = 1.95567621147336
$ws.Cells.Item(10, 3).Value = 0.203885460555739
$ws.Cells.Item(10, 4).Value = 0.1597301021452964
$ws.Cells.Item(10, 6).Value = 1.63264389447194
$ws.Cells.Item(10, 7).Value = 0.002463150854012553
$ws.Cells.Item(10, 10).Value = 0.1839748821900962
$ws.Cells.Item(10, 12).Value = 0.4099724289630018
$ws.Cells.Item(10, 14).Value = 1.381342435363322
$ws.Cells.Item(10, 15).Value = 4.077928409633785
$ws.Cells.Item(11, 2).Value = 2.038582378712363
$ws.Cells.Item(11, 3).Value = 0.2121458098633582
$ws.Cells.Item(11, 4).Value = 0.1604439688765353
$ws.Cells.Item(11, 6).Value = 1.635510867481003
$ws.Cells.Item(11, 7).Value = 0.002461116874493641
$ws.Cells.Item(11, 10).Value = 0.1829921080846297
$ws.Cells.Item(11, 12).Value = 0.4168724713372853
$ws.Cells.Item(11, 14).Value = 1.373224928281047
$ws.Cells.Item(11, 15).Value = 4.083744214090103
$ws.Cells.Item(12, 2).Value = 2.07003513517202
$ws.Cells.Item(12, 3).Value = 0.2152664272371396
$ws.Cells.Item(12, 4).Value = 0.1607209210622642
$ws.Cells.Item(12, 6).Value = 1.63673076070036
$ws.Cells.Item(12, 7).Value = 0.002460361536833711
$ws.Cells.Item(12, 10).Value = 0.1826300764778317
$ws.Cells.Item(12, 12).Value = 0.419506308510762
$ws.Cells.Item(12, 14).Value = 1.37023261676358
$ws.Cells.Item(12, 15).Value = 4.086306557752778
$ws.Cells.Item(13, 2).Value = 2.063258668552237
$ws.Cells.Item(13, 3).Value = 0.2145946779890551
$ws.Cells.Item(13, 4).Value = 0.1606609804643071
$ws.Cells.Item(13, 6).Value = 1.636462064196294
$ws.Cells.Item(13, 7).Value = 0.002460523551108364
$ws.Cells.Item(13, 10).Value = 0.1827075963434339
$ws.Cells.Item(13, 12).Value = 0.4189381359097126
$ws.Cells.Item(13, 14).Value = 1.370873435160306
$ws.Cells.Item(13, 15).Value = 4.085738692758667
$ws.Cells.Item(14, 2).Value = 2.041168863196333
$ws.Cells.Item(14, 3).Value = 0.2124026941017405
$ws.Cells.Item(14, 4).Value = 0.1604666213559582
$ws.Cells.Item(14, 6).Value = 1.635608538378065
$ws.Cells.Item(14, 7).Value = 0.002461054434434377
$ws.Cells.Item(14, 10).Value = 0.1829621206817116
$ws.Cells.Item(14, 12).Value = 0.4170887399951084
$ws.Cells.Item(14, 14).Value = 1.372977113489625
$ws.Cells.Item(14, 15).Value = 4.083947801968861
$ws.Cells.Item(15, 2).Value = 2.027645722445129
$ws.Cells.Item(15, 3).Value = 0.2110590725782799
$ws.Cells.Item(15, 4).Value = 0.1603484323949687
$ws.Cells.Item(15, 6).Value = 1.635103211418183
$ws.Cells.Item(15, 7).Value = 0.002461381552357612
$ws.Cells.Item(15, 10).Value = 0.1831193423352602
$ws.Cells.Item(15, 12).Value = 0.4159586536637363
$ws.Cells.Item(15, 14).Value = 1.374276304961008
$ws.Cells.Item(15, 15).Value = 4.082897727128255
$ws.Cells.Item(16, 2).Value = 1.950266303625313
$ws.Cells.Item(16, 3).Value = 0.2033446007539226
$ws.Cells.Item(16, 4).Value = 0.1596843789978308
$ws.Cells.Item(16, 6).Value = 1.632475312807898
$ws.Cells.Item(16, 7).Value = 0.002463285866246978
$ws.Cells.Item(16, 10).Value = 0.1840405258963891
$ws.Cells.Item(16, 12).Value = 0.4095244346148803
$ws.Cells.Item(16, 14).Value = 1.381884355081674
$ws.Cells.Item(16, 15).Value = 4.077598692042329
$ws.Cells.Item(17, 2).Value = 1.902901598150265
$ws.Cells.Item(17, 3).Value = 0.1985990158973152
$ws.Cells.Item(17, 4).Value = 0.1592888544853253
$ws.Cells.Item(17, 6).Value = 1.63110222529231
$ws.Cells.Item(17, 7).Value = 0.002464480690310195
$ws.Cells.Item(17, 10).Value = 0.184623682559014
$ws.Cells.Item(17, 12).Value = 0.405614730265313
$ws.Cells.Item(17, 14).Value = 1.386697031279134
$ws.Cells.Item(17, 15).Value = 4.074988736562716
$ws.Cells.Item(18, 2).Value = 1.875697815471426
$ws.Cells.Item(18, 3).Value = 0.1958647523401282
$ws.Cells.Item(18, 4).Value = 0.1590657289080823
$ws.Cells.Item(18, 6).Value = 1.630400296492326
$ws.Cells.Item(18, 7).Value = 0.002465177715278544
$ws.Cells.Item(18, 10).Value = 0.1849657313323476
$ws.Cells.Item(18, 12).Value = 0.4033798005471851
$ws.Cells.Item(18, 14).Value = 1.389518590791219
$ws.Cells.Item(18, 15).Value = 4.07372293660012
$ws.Cells.Item(19, 2).Value = 1.866493855597525
$ws.Cells.Item(19, 3).Value = 0.1949381703800555
$ws.Cells.Item(19, 4).Value = 0.1589909342965612
$ws.Cells.Item(19, 6).Value = 1.630177719973219
$ws.Cells.Item(19, 7).Value = 0.002465415400301893
$ws.Cells.Item(19, 10).Value = 0.1850826825640972
$ws.Cells.Item(19, 12).Value = 0.4026254720603788
$ws.Cells.Item(19, 14).Value = 1.390483100710661
$ws.Cells.Item(19, 15).Value = 4.0733347728588
$ws.Cells.Item(20, 2).Value = 1.907939612768359
$ws.Cells.Item(20, 3).Value = 0.199104682058703
$ws.Cells.Item(20, 4).Value = 0.1593305068304645
$ws.Cells.Item(20, 6).Value = 1.631239301938066
$ws.Cells.Item(20, 7).Value = 0.002464352485935414
$ws.Cells.Item(20, 10).Value = 0.1845609181606029
$ws.Cells.Item(20, 12).Value = 0.4060294949748311
$ws.Cells.Item(20, 14).Value = 1.386179183486355
$ws.Cells.Item(20, 15).Value = 4.075242207256082
$ws.Cells.Item(21, 2).Value = 2.047655613859661
$ws.Cells.Item(21, 3).Value = 0.2130467349372793
$ws.Cells.Item(21, 4).Value = 0.1605235298603915
$ws.Cells.Item(21, 6).Value = 1.635855596244795
$ws.Cells.Item(21, 7).Value = 0.002460898097883635
$ws.Cells.Item(21, 10).Value = 0.1828870860280745
$ws.Cells.Item(21, 12).Value = 0.4176313854330544
$ws.Cells.Item(21, 14).Value = 1.372356997388735
$ws.Cells.Item(21, 15).Value = 4.084464055914168
$ws.Cells.Item(22, 2).Value = 2.139305341601357
$ws.Cells.Item(22, 3).Value = 0.2221154799609053
$ws.Cells.Item(22, 4).Value = 0.1613418419048358
$ws.Cells.Item(22, 6).Value = 1.639655041236907
$ws.Cells.Item(22, 7).Value = 0.002458727190799504
$ws.Cells.Item(22, 10).Value = 0.1818521443644343
$ws.Cells.Item(22, 12).Value = 0.4253358721440605
$ws.Cells.Item(22, 14).Value = 1.363799091605763
$ws.Cells.Item(22, 15).Value = 4.092589748820302
$ws.Cells.Item(23, 2).Value = 2.09035984639246
$ws.Cells.Item(23, 3).Value = 0.2172793245773903
$ws.Cells.Item(23, 4).Value = 0.1609015761474097
$ws.Cells.Item(23, 6).Value = 1.637555600402123
$ws.Cells.Item(23, 7).Value = 0.00245987793141193
$ws.Cells.Item(23, 10).Value = 0.1823991160376721
$ws.Cells.Item(23, 12).Value = 0.4212127375883057
$ws.Cells.Item(23, 14).Value = 1.368323087161045
$ws.Cells.Item(23, 15).Value = 4.088060754563372
$ws.Cells.Item(24, 2).Value = 1.905661842697384
$ws.Cells.Item(24, 3).Value = 0.1988760889514083
$ws.Cells.Item(24, 4).Value = 0.1593116625126356
$ws.Cells.Item(24, 6).Value = 1.631177057099819
$ws.Cells.Item(24, 7).Value = 0.002464410415730082
$ws.Cells.Item(24, 10).Value = 0.1845892728120653
$ws.Cells.Item(24, 12).Value = 0.405841939931932
$ws.Cells.Item(24, 14).Value = 1.386413132118449
$ws.Cells.Item(24, 15).Value = 4.075126882104087
$ws.Cells.Item(25, 2).Value = 1.707937653241288
$ws.Cells.Item(25, 3).Value = 0.1788450113056967
$ws.Cells.Item(25, 4).Value = 0.1577639192807325
$ws.Cells.Item(25, 6).Value = 1.627654617411196
$ws.Cells.Item(25, 7).Value = 0.00246967418000262
$ws.Cells.Item(25, 10).Value = 0.1872053239408764
$ws.Cells.Item(25, 12).Value = 0.3897913716986636
$ws.Cells.Item(25, 14).Value = 1.407969288058069
$ws.Cells.Item(25, 15).Value = 4.070172748296528
